$d = $word.ActiveDocument

# Change 1: GPA 3.84 -> 3.86 ("." + "8" stays, "4" becomes "6")
$d.Content.Find.Execute(".84", $false, $false, $false, $false, $false, `
                         $true, 1, $false, ".86", 2)

# Change 2: append ", C/C++" after "Javascript, Java"
$d.Content.Find.Execute("Javascript, Java", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "Javascript, Java, C/C++", 2)
